# Update handback-status report timestamps / status after report regeneration.
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for rows 2 and 5
$overview.Range("G2").Value = "2016-08-22 14:16:28"
$overview.Range("G5").Value = "2016-08-22 14:16:28"

# zh-cn sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$zhcn.Range("E2").Value = "mt"
$zhcn.Range("E5").Value = "mt"
$zhcn.Range("H2").Value = "2016-08-22 14:16:23"
$zhcn.Range("H5").Value = "2016-08-22 14:16:23"
$zhcn.Range("K2").Value = "2016-08-22 14:16:40"
$zhcn.Range("K5").Value = "2016-08-22 14:16:40"

# de-de sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$dede.Range("E2").Value = "mt"
$dede.Range("E5").Value = "mt"
$dede.Range("H2").Value = "2016-08-22 14:16:28"
$dede.Range("H5").Value = "2016-08-22 14:16:28"
$dede.Range("K2").Value = "2016-08-22 14:16:47"
$dede.Range("K5").Value = "2016-08-22 14:16:47"
